$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "HET" sheet: add a second TFM_INS table (geothermal NCAP_BND limit) and
#    later rename the sheet itself to "PWR".
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HET")
$wsTra = $wb.Worksheets.Item("TRA")

# --- values / formulas -----------------------------------------------------
$ws.Range("B7").Value = "~TFM_INS"

$ws.Range("B8").Value = "TimeSlice"
$ws.Range("C8").Value = "LimType"
$ws.Range("D8").Value = "Attribute"
$ws.Range("E8").Value = "Year"
$ws.Range("F8").Formula = "=IF(Regions!C`$3<>`"`",Regions!C`$3,`"*`")"
$ws.Range("G8").Formula = "=IF(Regions!D`$3<>`"`",Regions!D`$3,`"*`")"
$ws.Range("H8").Value = "Pset_Set"
$ws.Range("I8").Value = "Pset_PN"
$ws.Range("J8").Value = "Cset_CN"

# row 9 / row 10 - PWR-GEO is entered before NCAP_BND so the new shared
# strings land in the same order as the target workbook (PWR-GEO, NCAP_BND)
$ws.Range("H9").Value = "PWR-GEO"
$ws.Range("D9").Value = "NCAP_BND"
$ws.Range("E9").Value = 2018
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

$ws.Range("H10").Value = "PWR-GEO"
$ws.Range("D10").Value = "NCAP_BND"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5

# --- formatting (mirror the existing table at row 2-4) ---------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null

$ws.Range("B3:I3").Copy() | Out-Null
$ws.Range("B8:I8").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Copy() | Out-Null
$ws.Range("J8").PasteSpecial(-4122) | Out-Null

# trailing marker cell (mirrors TRA!M2, the equivalent spot on another
# ~TFM_INS table in this workbook)
$wsTra.Range("M2").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- column width for the (now wider) Attribute column ---------------------
$ws.Columns.Item(4).ColumnWidth = 10

# --- view state --------------------------------------------------------------
$ws.Range("H17").Select()

# finally rename the sheet (after all "HET" lookups are done)
$ws.Name = "PWR"

# ---------------------------------------------------------------------------
# 2. "TRA" sheet: update the scrolled/selected view state.
# ---------------------------------------------------------------------------
$wsTra.Activate()
$wsTra.Range("B10:AI11").Select()
$excel.ActiveWindow.ScrollColumn = $wsTra.Range("V1").Column

# ---------------------------------------------------------------------------
# 3. Re-activate the PWR sheet (it is the active tab in the saved workbook).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("H17").Select()
